# Update "想去人数" (want-to-go count) values for several events.
# Same events appear on both the "展览" sheet and the "全部类型" sheet,
# but at different row numbers because "全部类型" interleaves rows from
# other categories.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1695
$ws1.Range("F9").Value = 67
$ws1.Range("F11").Value = 5971
$ws1.Range("F15").Value = 4653
$ws1.Range("F18").Value = 1154
$ws1.Range("F24").Value = 18
$ws1.Range("F25").Value = 3073
$ws1.Range("F26").Value = 132

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1695
$ws4.Range("F10").Value = 67
$ws4.Range("F12").Value = 5971
$ws4.Range("F16").Value = 4653
$ws4.Range("F19").Value = 1154
$ws4.Range("F25").Value = 18
$ws4.Range("F26").Value = 3073
$ws4.Range("F28").Value = 132
